$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This report ("localization-status.xlsx") tracks handoff/handback status for
# two source files (77bde751-... and d4dc8520-...) across the zh-cn and
# de-de locales. This edit records that the 77bde751 file has now been
# handed back and is in sync with en-US: its status flips from
# "Ready for handoff" to "Handed back: in sync with en-US" on every sheet,
# and the per-locale sheets gain the handback target/file links plus a
# handback timestamp.
# ---------------------------------------------------------------------------

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusHandedBack
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d50f05cd7d4681b034df6f614684e44eb2d83b36/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dba8594145e8f5fbfa7ce1a352c3b453a0e8bdb6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.zh-cn.xlf", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.zh-cn.xlf")
$zhcn.Range("H2").Value = "2016-03-20 06:29:21"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusHandedBack
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/d50f05cd7d4681b034df6f614684e44eb2d83b36/e2e/77bde751-9604-4978-951d-bf3e7caef7fd.md", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/20d73b9496eea8d892271f2bc2a0b76aa42aedcc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.de-de.xlf", "", "", "77bde751-9604-4978-951d-bf3e7caef7fd.1febe2f7767a8891db474e9d95c5b92db5df0285.de-de.xlf")
$dede.Range("H2").Value = "2016-03-20 06:29:27"
